$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFmt = "[`$-409]d/mmm/yyyy;@"
$amountFmt = "`"₹`"#,##0;`"₹`"\-#,##0"

# Row 249
$ws.Range("A249").Value = 44784
$ws.Range("A249").NumberFormat = $dateFmt
$ws.Range("B249").Value = "KA03MW7443"
$ws.Range("C249").Value = "JAZZ"
$ws.Range("D249").Value = "PMS"
$ws.Range("E249").Value = "WORK DONE DELIVERED"
$ws.Range("F249").Value = 3096
$ws.Range("F249").NumberFormat = $amountFmt
$ws.Range("G249").Value = "P PAY"

# Row 250
$ws.Range("A250").Value = 44784
$ws.Range("A250").NumberFormat = $dateFmt
$ws.Range("B250").Value = "KL07CJ4109"
$ws.Range("C250").Value = "I10"
$ws.Range("D250").Value = "PMS"
$ws.Range("E250").Value = "WORK DONE DELIVERED"
$ws.Range("F250").Value = 4950
$ws.Range("F250").NumberFormat = $amountFmt
$ws.Range("G250").Value = "G PAY"

# Row 251
$ws.Range("A251").Value = 44785
$ws.Range("A251").NumberFormat = $dateFmt
$ws.Range("B251").Value = "TN18U0010"
$ws.Range("C251").Value = "VERNA"
$ws.Range("D251").Value = "PMS"
$ws.Range("E251").Value = "WORK DONE DELIVERED"
$ws.Range("F251").Value = 4196
$ws.Range("F251").NumberFormat = $amountFmt

# Row 252
$ws.Range("A252").Value = 44785
$ws.Range("A252").NumberFormat = $dateFmt
$ws.Range("B252").Value = "KA01MG3419"
$ws.Range("C252").Value = "I20"
$ws.Range("D252").Value = "BODY SHOP"
$ws.Range("E252").Value = "WORK DONE"
$ws.Range("F252").Value = 26043
$ws.Range("F252").NumberFormat = $amountFmt
$ws.Range("G252").Value = "  INSURANCE"

# Row 253
$ws.Range("A253").Value = 44786
$ws.Range("A253").NumberFormat = $dateFmt
$ws.Range("B253").Value = "KA51MN5718"
$ws.Range("C253").Value = "VENTO"
$ws.Range("D253").Value = "RUNNING REPAIR"
$ws.Range("E253").Value = "WORK DONE DELIVERED"
$ws.Range("F253").Value = 6861
$ws.Range("F253").NumberFormat = $amountFmt
$ws.Range("G253").Value = "CARD"

# Row 254
$ws.Range("A254").Value = 44786
$ws.Range("A254").NumberFormat = $dateFmt
$ws.Range("B254").Value = "KA03MS0598"
$ws.Range("C254").Value = "A-STAR"
$ws.Range("D254").Value = "BODY SHOP"
$ws.Range("E254").Value = "WORK DONE DELIVERED"
$ws.Range("F254").Value = 5815
$ws.Range("F254").NumberFormat = $amountFmt
$ws.Range("G254").Value = "  INSURANCE"

# Row 255
$ws.Range("A255").Value = 44786
$ws.Range("A255").NumberFormat = $dateFmt
$ws.Range("B255").Value = "KA03AC3525"
$ws.Range("C255").Value = "DZIRE"
$ws.Range("D255").Value = "RUNNING REPAIR"
$ws.Range("E255").Value = "WORK DONE DELIVERED"
$ws.Range("F255").Value = 5479
$ws.Range("F255").NumberFormat = $amountFmt
$ws.Range("G255").Value = "P PAY"

# Row 256
$ws.Range("A256").Value = 44786
$ws.Range("A256").NumberFormat = $dateFmt
$ws.Range("B256").Value = "DL7CP8758"
$ws.Range("C256").Value = "I20"
$ws.Range("D256").Value = "RUNNING REPAIR"
$ws.Range("E256").Value = "WORK DONE DELIVERED"
$ws.Range("F256").Value = 14749
$ws.Range("F256").NumberFormat = $amountFmt

# Row 257
$ws.Range("A257").Value = 44786
$ws.Range("A257").NumberFormat = $dateFmt
$ws.Range("B257").Value = "KA03MS0598"
$ws.Range("C257").Value = "A-STAR"
$ws.Range("D257").Value = "PMS"
$ws.Range("E257").Value = "WORK DONE DELIVERED"
$ws.Range("F257").Value = 5002
$ws.Range("F257").NumberFormat = $amountFmt

# Row 258
$ws.Range("A258").Value = 44786
$ws.Range("A258").NumberFormat = $dateFmt
$ws.Range("B258").Value = "KA53MA4313"
$ws.Range("C258").Value = "RITZ"
$ws.Range("D258").Value = "MOULDING ROOF CHANGE"
$ws.Range("E258").Value = "WORK DONE DELIVERED"
$ws.Range("F258").Value = 600
$ws.Range("F258").NumberFormat = $amountFmt
$ws.Range("G258").Value = "CREDIT"

# Row 259
$ws.Range("A259").Value = 44786
$ws.Range("A259").NumberFormat = $dateFmt
$ws.Range("B259").Value = "KA03MZ4450"
$ws.Range("C259").Value = "BREZZA"
$ws.Range("D259").Value = "RUNNING REPAIR"
$ws.Range("E259").Value = "WORK DONE DELIVERED"
$ws.Range("F259").Value = 20459
$ws.Range("F259").NumberFormat = $amountFmt
$ws.Range("G259").Value = "CREDIT"

# Row 260
$ws.Range("A260").Value = 44788
$ws.Range("A260").NumberFormat = $dateFmt
$ws.Range("B260").Value = "KA51MK4755"
$ws.Range("C260").Value = "CRETA"
$ws.Range("D260").Value = "PMS"
$ws.Range("E260").Value = "WORK DONE DELIVERED"
$ws.Range("F260").Value = 6433
$ws.Range("F260").NumberFormat = $amountFmt
$ws.Range("G260").Value = "CARD"

# Row 261
$ws.Range("A261").Value = 44788
$ws.Range("A261").NumberFormat = $dateFmt
$ws.Range("B261").Value = "MH43AT5184"
$ws.Range("C261").Value = "I20"
$ws.Range("D261").Value = "PMS"
$ws.Range("E261").Value = "WORK DONE DELIVERED"
$ws.Range("F261").Value = 10277
$ws.Range("F261").NumberFormat = $amountFmt

# Row 262
$ws.Range("A262").Value = 44788
$ws.Range("A262").NumberFormat = $dateFmt
$ws.Range("B262").Value = "AP16AV0759"
$ws.Range("C262").Value = "M OMNI"
$ws.Range("D262").Value = "CLUTCH PROBLEM"
$ws.Range("E262").Value = "WORK IN PROGRESS"

# Row 263
$ws.Range("A263").Value = 44788
$ws.Range("A263").NumberFormat = $dateFmt
$ws.Range("B263").Value = "KA01MR6461"
$ws.Range("C263").Value = "120 ASTA"
$ws.Range("D263").Value = "CLUTCH PROBLEM"
$ws.Range("E263").Value = "WORK IN PROGRESS"

# Row 264
$ws.Range("A264").Value = 44788
$ws.Range("A264").NumberFormat = $dateFmt
$ws.Range("B264").Value = "KA03MT2662"
$ws.Range("C264").Value = "POLO"
$ws.Range("D264").Value = "RUNNING REPAIR"
$ws.Range("E264").Value = "WORK IN PROGRESS"

# Update the active view: scroll position + selection to match the final edit state
$win = $excel.Windows.Item(1)
$win.ScrollRow = 241
$win.ScrollColumn = 1
$ws.Range("G261").Select()
